$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "61.034.89"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.07%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.881.84"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -1.49%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "587.93"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.42%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "138.78"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -5.29%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.491"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -3.29%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "6.87"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.67%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.137"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -4.74%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.428"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -3.09%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000217"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -3.95%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "32.20"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -4.15%  "
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.33%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.361.67"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -1.48%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "60.985.99"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.09%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.882.85"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -1.45%  "
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -3.36%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "424.92"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -1.71%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.28"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.78%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.652"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -4.04%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.90"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -2.66%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "79.76"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -2.14%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.43"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -4.55%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.09%  "
$ws.Range("B26").NumberFormat = "@"
$ws.Range("B26").Value = "InternetComputer(DFINITY)"
$ws.Range("C26").NumberFormat = "@"
$ws.Range("C26").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.39"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -4.10%  "
$ws.Range("B27").NumberFormat = "@"
$ws.Range("B27").Value = "Fetch.AI"
$ws.Range("C27").NumberFormat = "@"
$ws.Range("C27").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.04"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -7.09%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -3.38%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.06"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -9.18%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.63"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -5.37%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.05%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "25.67"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -3.73%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -4.51%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0₃0842"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -1.68%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.968"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -4.45%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.42"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -3.90%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -2.31%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.79"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -7.01%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.89"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -4.71%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.31"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -3.15%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.115"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -5.25%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "38.22"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -4.47%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -6.73%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.657.73"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -1.61%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "132.15"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -1.26%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0328"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -4.35%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "345.99"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -9.16%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -4.09%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "22.22"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -6.74%  "
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -4.35%  "
